$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 78
$ws.Cells.Item(11, 9).Value = 78
$ws.Cells.Item(11, 11).Value = 78
$ws.Cells.Item(11, 13).Value = 62
$ws.Cells.Item(38, 8).Value = 261.9091
$ws.Cells.Item(38, 9).Value = 296.77777
$ws.Cells.Item(38, 10).Value = 105
$ws.Cells.Item(38, 11).Value = 890.33331
$ws.Cells.Item(38, 12).Value = 315
$ws.Cells.Item(38, 13).Value = -518.33331
$ws.Cells.Item(38, 14).Value = -1059
$ws.Cells.Item(39, 8).Value = 304.2
$ws.Cells.Item(39, 9).Value = 130.25
$ws.Cells.Item(39, 11).Value = 390.75
$ws.Cells.Item(39, 13).Value = -94.75
$ws.Cells.Item(62, 8).Value = 2428
$ws.Cells.Item(62, 9).Value = 1142
$ws.Cells.Item(62, 11).Value = 1142
$ws.Cells.Item(62, 13).Value = -518
$ws.Cells.Item(65, 8).Value = 2428
$ws.Cells.Item(65, 9).Value = 1142
$ws.Cells.Item(65, 11).Value = 5710
$ws.Cells.Item(65, 13).Value = -2590
$ws.Cells.Item(92, 8).Value = 5291271.5
$ws.Cells.Item(92, 9).Value = 290.55554
$ws.Cells.Item(92, 10).Value = 37037156
$ws.Cells.Item(92, 11).Value = 290.55554
$ws.Cells.Item(92, 12).Value = 37037156
$ws.Cells.Item(92, 13).Value = 957.4444599999999
$ws.Cells.Item(92, 14).Value = -37039652
$ws.Cells.Item(111, 8).Value = 4654
$ws.Cells.Item(111, 9).Value = 3461.25
$ws.Cells.Item(111, 11).Value = 10383.75
$ws.Cells.Item(111, 13).Value = -7316.75
$ws.Cells.Item(132, 8).Value = 20158.893
$ws.Cells.Item(132, 9).Value = 1429.6818
$ws.Cells.Item(132, 11).Value = 4289.0454
$ws.Cells.Item(132, 13).Value = -1759.0454
$ws.Cells.Item(137, 8).Value = 6338.9565
$ws.Cells.Item(137, 9).Value = 6753.615
$ws.Cells.Item(137, 11).Value = 20260.845
$ws.Cells.Item(137, 13).Value = -17710.845
$ws.Cells.Item(138, 8).Value = 3248.6191
$ws.Cells.Item(138, 9).Value = 1828
$ws.Cells.Item(138, 10).Value = 3636.0605
$ws.Cells.Item(138, 11).Value = 5484
$ws.Cells.Item(138, 12).Value = 10908.1815
$ws.Cells.Item(138, 13).Value = -344
$ws.Cells.Item(138, 14).Value = -21188.1815

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 22187.666
$ws.Cells.Item(2, 9).Value = 26671.357
$ws.Cells.Item(2, 11).Value = 26671.357
$ws.Cells.Item(2, 13).Value = -26558.357
$ws.Cells.Item(32, 8).Value = 2247.1633
$ws.Cells.Item(32, 9).Value = 2204.0977
$ws.Cells.Item(32, 11).Value = 2204.0977
$ws.Cells.Item(32, 13).Value = -1917.0977
$ws.Cells.Item(45, 8).Value = 7790.75
$ws.Cells.Item(45, 10).Value = 10749.5
$ws.Cells.Item(45, 12).Value = 10749.5
$ws.Cells.Item(45, 14).Value = -11503.5
$ws.Cells.Item(102, 8).Value = 30312888
$ws.Cells.Item(102, 10).Value = 111129180
$ws.Cells.Item(102, 12).Value = 111129180
$ws.Cells.Item(102, 14).Value = -111132424
$ws.Cells.Item(116, 8).Value = 22187.666
$ws.Cells.Item(116, 9).Value = 26671.357
$ws.Cells.Item(116, 11).Value = 26671.357
$ws.Cells.Item(116, 13).Value = -24377.357
$ws.Cells.Item(132, 8).Value = 1588.2273
$ws.Cells.Item(132, 9).Value = 1445.9445
$ws.Cells.Item(132, 11).Value = 4337.833500000001
$ws.Cells.Item(132, 13).Value = -1807.833500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 22187.666
$ws.Cells.Item(3, 9).Value = 26671.357
$ws.Cells.Item(3, 11).Value = 26671.357
$ws.Cells.Item(3, 13).Value = -26557.357
$ws.Cells.Item(22, 8).Value = 12924.5
$ws.Cells.Item(22, 9).Value = 50150.5
$ws.Cells.Item(22, 10).Value = 515.8333
$ws.Cells.Item(22, 11).Value = 50150.5
$ws.Cells.Item(22, 12).Value = 515.8333
$ws.Cells.Item(22, 13).Value = -49977.5
$ws.Cells.Item(22, 14).Value = -861.8333
$ws.Cells.Item(94, 8).Value = 3973165.2
$ws.Cells.Item(94, 9).Value = 4343.5713
$ws.Cells.Item(94, 10).Value = 15879630
$ws.Cells.Item(94, 11).Value = 4343.5713
$ws.Cells.Item(94, 12).Value = 15879630
$ws.Cells.Item(94, 13).Value = -3892.5713
$ws.Cells.Item(94, 14).Value = -15880532
$ws.Cells.Item(134, 8).Value = 1598.32
$ws.Cells.Item(134, 9).Value = 1636.8572
$ws.Cells.Item(134, 11).Value = 4910.571599999999
$ws.Cells.Item(134, 13).Value = -2375.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3306.1875
$ws.Cells.Item(16, 9).Value = 2116.5833
$ws.Cells.Item(16, 11).Value = 2116.5833
$ws.Cells.Item(16, 13).Value = -1829.5833
$ws.Cells.Item(31, 8).Value = 3754.0667
$ws.Cells.Item(31, 9).Value = 2979.111
$ws.Cells.Item(31, 11).Value = 2979.111
$ws.Cells.Item(31, 13).Value = -2684.111
$ws.Cells.Item(34, 8).Value = 3754.0667
$ws.Cells.Item(34, 9).Value = 2979.111
$ws.Cells.Item(34, 11).Value = 2979.111
$ws.Cells.Item(34, 13).Value = -2777.111
$ws.Cells.Item(113, 8).Value = 3306.1875
$ws.Cells.Item(113, 9).Value = 2116.5833
$ws.Cells.Item(113, 11).Value = 2116.5833
$ws.Cells.Item(113, 13).Value = 53.41670000000022
$ws.Cells.Item(134, 8).Value = 3288.7
$ws.Cells.Item(134, 9).Value = 3338.4644
$ws.Cells.Item(134, 10).Value = 2592
$ws.Cells.Item(134, 11).Value = 10015.3932
$ws.Cells.Item(134, 12).Value = 7776
$ws.Cells.Item(134, 13).Value = -7480.393199999999
$ws.Cells.Item(134, 14).Value = -12846

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 6353.2666
$ws.Cells.Item(80, 10).Value = 6422.1113
$ws.Cells.Item(80, 12).Value = 19266.3339
$ws.Cells.Item(80, 14).Value = -21138.3339
$ws.Cells.Item(83, 8).Value = 6353.2666
$ws.Cells.Item(83, 10).Value = 6422.1113
$ws.Cells.Item(83, 12).Value = 57799.00169999999
$ws.Cells.Item(83, 14).Value = -67159.00169999999
$ws.Cells.Item(88, 8).Value = 4714.143
$ws.Cells.Item(88, 9).Value = 3833.3333
$ws.Cells.Item(88, 10).Value = 5374.75
$ws.Cells.Item(88, 11).Value = 11499.9999
$ws.Cells.Item(88, 12).Value = 16124.25
$ws.Cells.Item(88, 13).Value = -11071.9999
$ws.Cells.Item(88, 14).Value = -16980.25
$ws.Cells.Item(91, 8).Value = 4714.143
$ws.Cells.Item(91, 9).Value = 3833.3333
$ws.Cells.Item(91, 10).Value = 5374.75
$ws.Cells.Item(91, 11).Value = 11499.9999
$ws.Cells.Item(91, 12).Value = 16124.25
$ws.Cells.Item(91, 13).Value = -10017.9999
$ws.Cells.Item(91, 14).Value = -19088.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 5882422.5
$ws.Cells.Item(2, 9).Value = 63.666668
$ws.Cells.Item(2, 10).Value = 12500077
$ws.Cells.Item(2, 11).Value = 63.666668
$ws.Cells.Item(2, 12).Value = 12500077
$ws.Cells.Item(2, 13).Value = 49.333332
$ws.Cells.Item(2, 14).Value = -12500303
$ws.Cells.Item(102, 8).Value = 3551
$ws.Cells.Item(102, 9).Value = 3586.8
$ws.Cells.Item(102, 11).Value = 3586.8
$ws.Cells.Item(102, 13).Value = -1964.8
$ws.Cells.Item(113, 8).Value = 22582.916
$ws.Cells.Item(113, 9).Value = 12199
$ws.Cells.Item(113, 11).Value = 12199
$ws.Cells.Item(113, 13).Value = -10029
$ws.Cells.Item(122, 8).Value = 4518.2085
$ws.Cells.Item(122, 9).Value = 3759.842
$ws.Cells.Item(122, 11).Value = 11279.526
$ws.Cells.Item(122, 13).Value = -8829.526

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 473.73914
$ws.Cells.Item(55, 9).Value = 434.33334
$ws.Cells.Item(55, 10).Value = 547.625
$ws.Cells.Item(55, 11).Value = 434.33334
$ws.Cells.Item(55, 12).Value = 547.625
$ws.Cells.Item(55, 13).Value = -261.33334
$ws.Cells.Item(55, 14).Value = -893.625
$ws.Cells.Item(122, 8).Value = 4250
$ws.Cells.Item(122, 9).Value = 4250
$ws.Cells.Item(122, 11).Value = 12750
$ws.Cells.Item(122, 13).Value = -10300

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 812.93335
$ws.Cells.Item(107, 10).Value = 789.6
$ws.Cells.Item(107, 12).Value = 2368.8
$ws.Cells.Item(107, 14).Value = -6208.8
$ws.Cells.Item(126, 8).Value = 3283
$ws.Cells.Item(126, 9).Value = 2609.4
$ws.Cells.Item(126, 10).Value = 4125
$ws.Cells.Item(126, 11).Value = 7828.200000000001
$ws.Cells.Item(126, 12).Value = 12375
$ws.Cells.Item(126, 13).Value = -5358.200000000001
$ws.Cells.Item(126, 14).Value = -17315
$ws.Cells.Item(139, 8).Value = 45999.2
$ws.Cells.Item(139, 10).Value = 40000
$ws.Cells.Item(139, 14).Value = -50280

Write-Output "Applied 184 cell updates"